$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell 'D2' '42.502.69'
Set-TextCell 'E2' '  -0.22%  '
Set-TextCell 'D3' '2.331.39'
Set-TextCell 'E3' '  -1.32%  '
Set-TextCell 'E4' '  +0.19%  '
Set-TextCell 'D5' '311.00'
Set-TextCell 'E5' '  -5.07%  '
Set-TextCell 'D6' '106.62'
Set-TextCell 'E6' '  +5.88%  '
Set-TextCell 'D7' '0.630'
Set-TextCell 'E7' '  -1.10%  '
Set-TextCell 'E8' '  +0.02%  '
Set-TextCell 'D9' '0.607'
Set-TextCell 'E9' '  -2.16%  '
Set-TextCell 'D10' '40.03'
Set-TextCell 'E10' '  -0.44%  '
Set-TextCell 'D11' '0.0920'
Set-TextCell 'E11' '  -0.16%  '
Set-TextCell 'D12' '8.40'
Set-TextCell 'E12' '  -0.17%  '
Set-TextCell 'B13' 'TRON'
Set-TextCell 'C13' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 'D13' '0.106'
Set-TextCell 'E13' '  +0.76%  '
Set-TextCell 'B14' 'Polygon'
Set-TextCell 'C14' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 'D14' '0.985'
Set-TextCell 'E14' '  -3.00%  '
Set-TextCell 'D15' '15.68'
Set-TextCell 'E15' '  -3.70%  '
Set-TextCell 'D16' '2.687.85'
Set-TextCell 'E16' '  -1.21%  '
Set-TextCell 'D17' '2.328.74'
Set-TextCell 'E17' '  -1.42%  '
Set-TextCell 'D18' '42.514.05'
Set-TextCell 'E18' '  -0.13%  '
Set-TextCell 'D19' '7.52'
Set-TextCell 'E19' '  -1.47%  '
Set-TextCell 'D20' '0.0000105'
Set-TextCell 'E20' '  -0.98%  '
Set-TextCell 'D21' '75.78'
Set-TextCell 'E21' '  +0.77%  '
Set-TextCell 'D22' '3.51'
Set-TextCell 'E22' '  -4.85%  '
Set-TextCell 'D23' '265.45'
Set-TextCell 'E23' '  -3.80%  '
Set-TextCell 'D24' '2.30'
Set-TextCell 'E24' '  -0.27%  '
Set-TextCell 'D25' '9.37'
Set-TextCell 'E25' '  -4.11%  '
Set-TextCell 'E26' '  +0.41%  '
Set-TextCell 'D27' '11.16'
Set-TextCell 'E27' '  -2.79%  '
Set-TextCell 'D28' '23.25'
Set-TextCell 'E28' '  -2.12%  '
Set-TextCell 'E29' '  +1.80%  '
Set-TextCell 'D30' '35.85'
Set-TextCell 'E30' '  +0.89%  '
Set-TextCell 'D31' '165.96'
Set-TextCell 'E31' '  -4.54%  '
Set-TextCell 'D32' '0.0897'
Set-TextCell 'E32' '  -0.34%  '
Set-TextCell 'B33' 'WEMIXToken'
Set-TextCell 'C33' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell 'D33' '2.92'
Set-TextCell 'E33' '  -5.87%  '
Set-TextCell 'B34' 'Filecoin'
Set-TextCell 'C34' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D34' '5.98'
Set-TextCell 'E34' '  +0.50%  '
Set-TextCell 'D35' '0.120'
Set-TextCell 'E35' '  +15.11%  '
Set-TextCell 'D36' '0.130'
Set-TextCell 'E36' '  -1.88%  '
Set-TextCell 'D37' '4.57'
Set-TextCell 'E37' '  -0.87%  '
Set-TextCell 'D38' '0.0355'
Set-TextCell 'E38' '  -0.83%  '
Set-TextCell 'D39' '3.76'
Set-TextCell 'E39' '  -3.32%  '
Set-TextCell 'D40' '2.63'
Set-TextCell 'E40' '  -8.43%  '
Set-TextCell 'D41' '104.12'
Set-TextCell 'E41' '  +15.51%  '
Set-TextCell 'B42' 'Algorand'
Set-TextCell 'C42' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D42' '0.235'
Set-TextCell 'E42' '  +3.05%  '
Set-TextCell 'B43' 'ARBITRUM'
Set-TextCell 'C43' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 'D43' '1.47'
Set-TextCell 'E43' '  -2.99%  '
Set-TextCell 'D44' '70.82'
Set-TextCell 'E44' '  +3.06%  '
Set-TextCell 'E45' '  +0.28%  '
Set-TextCell 'D46' '12.18'
Set-TextCell 'E46' '  +2.20%  '
Set-TextCell 'D47' '112.07'
Set-TextCell 'E47' '  -2.78%  '
Set-TextCell 'B48' 'THORChain'
Set-TextCell 'C48' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextCell 'D48' '5.45'
Set-TextCell 'E48' '  -0.04%  '
Set-TextCell 'B49' 'ordi'
Set-TextCell 'C49' 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextCell 'D49' '76.86'
Set-TextCell 'E49' '  +12.88%  '
Set-TextCell 'D50' '9.00'
Set-TextCell 'E50' '  -0.23%  '
Set-TextCell 'D51' '1.27'
Set-TextCell 'E51' '  +0.24%  '
